$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (shared-string rich-text runs) ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Step 1: fix up cells whose type (text vs numeric) changes, by cloning format+type from a stable template cell ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C16").Copy($ws.Range("D20"))
$ws.Range("E16").Copy($ws.Range("E20"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C16").Copy($ws.Range("D30"))
$ws.Range("E16").Copy($ws.Range("E30"))

# --- Step 2: write the final numeric/text values for every changed cell ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 54
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = -10
$ws.Range("L16").Value = 54.285714285714
$ws.Range("M16").Value = 58.823529411764
$ws.Range("N16").Value = -77.5
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 53
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = 17.777777777777
$ws.Range("L17").Value = 43.243243243243
$ws.Range("M17").Value = 1.923076923076
$ws.Range("N17").Value = -40.449438202247
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = -52.631578947368
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 63
$ws.Range("K18").Value = -30.158730158730
$ws.Range("L18").Value = 2.325581395348
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = -77.319587628866
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -31.25
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 276
$ws.Range("J19").Value = 244
$ws.Range("K19").Value = 13.114754098360
$ws.Range("L19").Value = 79.220779220779
$ws.Range("M19").Value = 22.666666666666
$ws.Range("N19").Value = -15.596330275229
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 91.666666666666
$ws.Range("M20").Value = 35.294117647058
$ws.Range("N20").Value = -88.942307692307
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -35.483870967741
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = 2.222222222222
$ws.Range("I21").Value = 452
$ws.Range("J21").Value = 441
$ws.Range("K21").Value = 2.494331065759
$ws.Range("L21").Value = 58.596491228070
$ws.Range("M21").Value = 21.505376344086
$ws.Range("N21").Value = -57.717492984097
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("M22").Value = 300
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 25
$ws.Range("K23").Value = -28
$ws.Range("L23").Value = -41.935483870967
$ws.Range("M23").Value = -5.263157894736
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -15.789473684210
$ws.Range("F24").Value = 61
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = 8.928571428571
$ws.Range("I24").Value = 261
$ws.Range("J24").Value = 246
$ws.Range("K24").Value = 6.097560975609
$ws.Range("L24").Value = 39.572192513369
$ws.Range("M24").Value = -15.806451612903
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -51.020408163265
$ws.Range("I25").Value = 130
$ws.Range("J25").Value = 127
$ws.Range("K25").Value = 2.362204724409
$ws.Range("L25").Value = 44.444444444444
$ws.Range("M25").Value = -5.109489051094
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = 12.5
